# Decrement column E (剩余) by 1 for every data row (2-99),
# except row 36 which keeps its original value (not part of this update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # Column E is the 5th column
    $cell.Value2 = $cell.Value2 - 1
}
